$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.213166144200627
$ws.Range("C2").Value = 0.5297805642633229
$ws.Range("J2").Value = 0.01567398119122257
$ws.Range("P2").Value = 0.1536050156739812
$ws.Range("S2").Value = 0.0877742946708464
$ws.Range("B3").Value = 0.005617977528089887
$ws.Range("C3").Value = 0.05056179775280899
$ws.Range("J3").Value = 0.02808988764044944
$ws.Range("P3").Value = 0.7359550561797753
$ws.Range("S3").Value = 0.1797752808988764
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7804878048780488
$ws.Range("S4").Value = 0.1951219512195122
$ws.Range("B6").Value = 0.0821917808219178
$ws.Range("D6").Value = 0.0091324200913242
$ws.Range("F6").Value = 0.0410958904109589
$ws.Range("J6").Value = 0.273972602739726
$ws.Range("O6").Value = 0.0319634703196347
$ws.Range("Q6").Value = 0.1643835616438356
$ws.Range("R6").Value = 0.091324200913242
$ws.Range("S6").Value = 0.3059360730593607
$ws.Range("B7").Value = 0.08441558441558442
$ws.Range("D7").Value = 0.01298701298701299
$ws.Range("F7").Value = 0.08441558441558442
$ws.Range("J7").Value = 0.1038961038961039
$ws.Range("O7").Value = 0.01298701298701299
$ws.Range("Q7").Value = 0.2142857142857143
$ws.Range("R7").Value = 0.05844155844155844
$ws.Range("S7").Value = 0.4285714285714285
$ws.Range("B8").Value = 0.1234567901234568
$ws.Range("D8").Value = 0.02962962962962963
$ws.Range("E8").Value = 0.004938271604938272
$ws.Range("F8").Value = 0.07407407407407407
$ws.Range("J8").Value = 0.1061728395061728
$ws.Range("O8").Value = 0.03209876543209877
$ws.Range("Q8").Value = 0.2172839506172839
$ws.Range("R8").Value = 0.08641975308641975
$ws.Range("S8").Value = 0.325925925925926
$ws.Range("B9").Value = 0.130718954248366
$ws.Range("D9").Value = 0.006535947712418301
$ws.Range("F9").Value = 0.08496732026143791
$ws.Range("J9").Value = 0.1633986928104575
$ws.Range("O9").Value = 0.006535947712418301
$ws.Range("Q9").Value = 0.1699346405228758
$ws.Range("R9").Value = 0.130718954248366
$ws.Range("S9").Value = 0.3071895424836601
$ws.Range("B10").Value = 0.1291208791208791
$ws.Range("D10").Value = 0.02380952380952381
$ws.Range("E10").Value = 0.003663003663003663
$ws.Range("J10").Value = 0.1291208791208791
$ws.Range("O10").Value = 0.02197802197802198
$ws.Range("Q10").Value = 0.2032967032967033
$ws.Range("R10").Value = 0.1034798534798535
$ws.Range("S10").Value = 0.3141025641025641
$ws.Range("G11").Value = 0.140495867768595
$ws.Range("J11").Value = 0.09504132231404959
$ws.Range("K11").Value = 0.2190082644628099
$ws.Range("L11").Value = 0.5330578512396694
$ws.Range("S11").Value = 0.01239669421487603
$ws.Range("G12").Value = 0.7727272727272727
$ws.Range("J12").Value = 0.1515151515151515
$ws.Range("K12").Value = 0.007575757575757576
$ws.Range("L12").Value = 0.03787878787878788
$ws.Range("G13").Value = 0.6571428571428571
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.02666666666666667
$ws.Range("H15").Value = 0.1466666666666667
$ws.Range("I15").Value = 0.03555555555555556
$ws.Range("J15").Value = 0.3466666666666667
$ws.Range("K15").Value = 0.07111111111111111
$ws.Range("M15").Value = 0.008888888888888889
$ws.Range("O15").Value = 0.1022222222222222
$ws.Range("S15").Value = 0.2622222222222222
$ws.Range("F16").Value = 0.00966183574879227
$ws.Range("H16").Value = 0.1980676328502415
$ws.Range("I16").Value = 0.05797101449275362
$ws.Range("J16").Value = 0.3719806763285024
$ws.Range("K16").Value = 0.09178743961352658
$ws.Range("M16").Value = 0.02415458937198068
$ws.Range("N16").Value = 0.004830917874396135
$ws.Range("O16").Value = 0.07729468599033816
$ws.Range("S16").Value = 0.1642512077294686
$ws.Range("F17").Value = 0.02463054187192118
$ws.Range("H17").Value = 0.1896551724137931
$ws.Range("I17").Value = 0.08866995073891626
$ws.Range("J17").Value = 0.4187192118226601
$ws.Range("K17").Value = 0.103448275862069
$ws.Range("M17").Value = 0.01477832512315271
$ws.Range("N17").Value = 0.002463054187192118
$ws.Range("S17").Value = 0.08620689655172414
$ws.Range("F18").Value = 0.03045685279187817
$ws.Range("H18").Value = 0.2284263959390863
$ws.Range("I18").Value = 0.116751269035533
$ws.Range("J18").Value = 0.3248730964467005
$ws.Range("K18").Value = 0.1015228426395939
$ws.Range("M18").Value = 0.01015228426395939
$ws.Range("O18").Value = 0.08629441624365482
$ws.Range("S18").Value = 0.1015228426395939
$ws.Range("F19").Value = 0.02074688796680498
$ws.Range("H19").Value = 0.2147302904564315
$ws.Range("I19").Value = 0.07676348547717843
$ws.Range("J19").Value = 0.3921161825726141
$ws.Range("K19").Value = 0.09647302904564316
$ws.Range("M19").Value = 0.02074688796680498
$ws.Range("N19").Value = 0.003112033195020747
$ws.Range("O19").Value = 0.07676348547717843
$ws.Range("S19").Value = 0.09854771784232365
